# Applies the edits described by the commit:
#  1. Add comment "Not a tool" anchored on the "p" of "poles"
#     (Others place two well-staked poles...)
#  2. Add comment "Tag as tool" anchored on the "r" of "ruler"
#     (And underneath they place the above-mentioned ruler.)
#  3. "a similar tool with measured intervals" -> "a similar tool divided with measured intervals"
#  4. Add comment "Part of the canon, not a tool" anchored on the "t" of "touch hole"
#     (...seek its exit through the touch hole which is spiked...)
#  5. "you will make a nut in the second touchhole" -> "you will make a <tl>nut</tl> in the second touchhole"
#
# NOTE: in this runtime, writing properties (Author/Initial) through the
# Comment object returned directly by Comments.Add() is unreliable once
# more than one comment exists -- the write silently lands on the wrong
# comment. Re-fetching the comment via $d.Comments.Item(<1-based index>)
# right before setting properties works reliably, so that pattern is used
# throughout below.

$d = $word.ActiveDocument

function Add-CommentOnSubstring($searchText, $needle, $commentText) {
    # Re-finds searchText fresh (so earlier edits/comments do not invalidate
    # offsets), locates the first occurrence of $needle inside the matched
    # text, and anchors a new comment on that single character -- mirroring
    # the <w:commentRangeStart/.../w:commentRangeEnd/commentReference> split
    # seen in the target diff.
    $r = $d.Content
    $found = $r.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find: $searchText"
    }
    $matchStart = $r.Start
    $offset = $r.Text.IndexOf($needle)
    if ($offset -lt 0) {
        throw "Needle '$needle' not found inside '$searchText'"
    }
    $charStart = $matchStart + $offset
    $charEnd = $charStart + 1
    $charRange = $d.Range($charStart, $charEnd)

    $d.Comments.Add($charRange, $commentText) | Out-Null
    $newIndex = $d.Comments.Count
    $d.Comments.Item($newIndex).Author = "Margot Lyautey"
    $d.Comments.Item($newIndex).Initial = "ML"
}

# --- Comment 0: "poles" -> split off the leading "p" ----------------------
Add-CommentOnSubstring "Others place two well-staked poles" "poles" "Not a tool"

# --- Comment 1: "ruler" -> split off the leading "r" -----------------------
Add-CommentOnSubstring "the above-mentioned ruler. Then when want to fire" "ruler" "Tag as tool"

# --- Text change: "a similar tool with measured intervals" ----------------
$d.Content.Find.Execute("a similar tool with measured intervals", $true, $false, $false, $false, $false, $true, 1, $false, "a similar tool divided with measured intervals", 2) | Out-Null

# --- Comment 2: "touch hole" -> split off the leading "t" ------------------
Add-CommentOnSubstring "seek its exit through the touch hole which is spiked" "touch" "Part of the canon, not a tool"

# --- Text change: "make a nut" -> "make a <tl>nut</tl>" --------------------
$d.Content.Find.Execute("you will make a nut in the second touchhole", $true, $false, $false, $false, $false, $true, 1, $false, "you will make a <tl>nut</tl> in the second touchhole", 2) | Out-Null
